# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Sat Jul 27 11:58:00 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.191.52"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "3.278.42"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.49%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  +4.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("D12").Value = "3.846.71"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.137"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.48%  "
$ws.Range("D15").Value = "68.199.56"
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("E16").Value = "  +1.95%  "
$ws.Range("D17").Value = "3.273.99"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("E19").Value = "  +2.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "382.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.41%  "
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("E26").Value = "  +7.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.995"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("E29").Value = "  +3.24%  "
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.01%  "
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("E35").Value = "  +2.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.840"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.21%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.06%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("E44").Value = "  +2.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.99%  "
$ws.Range("D46").Value = "2.628.60"
$ws.Range("E46").Value = "  -5.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "342.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.01%  "
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.37%  "
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("E51").Value = "  -0.20%  "
